# Croatia 3NL - Atualizacao de bases das ligas, do dia: 16-06-2024 as 07:16
#
# The underlying source data had several fixture rows whose match data (everything
# except the running id in column A, the Div in column C and the Date in column D)
# had been shuffled between rows. This script restores the correct pairing by
# swapping (or rotating) the affected columns back between the impacted rows.
#
# NOTE: this sandboxed PowerShell engine does not create separate variable scopes
# for functions, so every helper below is careful to use unique variable names
# that cannot collide with any other variable used anywhere else in the script.

function Get-MatchRowValues {
    param($rowValuesSheet, $rowValuesRow, $rowValuesCols)
    $rowValuesResult = @{}
    foreach ($rowValuesColName in $rowValuesCols) {
        $rowValuesResult[$rowValuesColName] = $rowValuesSheet.Range($rowValuesColName + $rowValuesRow).Value()
    }
    return $rowValuesResult
}

function Set-MatchRowValues {
    param($setRowSheet, $setRowRow, $setRowCols, $setRowVals)
    foreach ($setRowColName in $setRowCols) {
        $setRowSheet.Range($setRowColName + $setRowRow).Value = $setRowVals[$setRowColName]
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns that carry the per-match data (everything except A=id, C=Div, D=Date)
$matchCols = @("B","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD")

# --- Pair swap: rows 7 and 10 ---
$saved007 = Get-MatchRowValues $ws 7 $matchCols
$saved010 = Get-MatchRowValues $ws 10 $matchCols
Set-MatchRowValues $ws 7 $matchCols $saved010
Set-MatchRowValues $ws 10 $matchCols $saved007

# --- Pair swap: rows 23 and 24 ---
$saved023 = Get-MatchRowValues $ws 23 $matchCols
$saved024 = Get-MatchRowValues $ws 24 $matchCols
Set-MatchRowValues $ws 23 $matchCols $saved024
Set-MatchRowValues $ws 24 $matchCols $saved023

# --- Pair swap: rows 29 and 30 ---
$saved029 = Get-MatchRowValues $ws 29 $matchCols
$saved030 = Get-MatchRowValues $ws 30 $matchCols
Set-MatchRowValues $ws 29 $matchCols $saved030
Set-MatchRowValues $ws 30 $matchCols $saved029

# --- Pair swap: rows 72 and 73 ---
$saved072 = Get-MatchRowValues $ws 72 $matchCols
$saved073 = Get-MatchRowValues $ws 73 $matchCols
Set-MatchRowValues $ws 72 $matchCols $saved073
Set-MatchRowValues $ws 73 $matchCols $saved072

# --- Pair swap: rows 127 and 129 ---
$saved127 = Get-MatchRowValues $ws 127 $matchCols
$saved129 = Get-MatchRowValues $ws 129 $matchCols
Set-MatchRowValues $ws 127 $matchCols $saved129
Set-MatchRowValues $ws 129 $matchCols $saved127

# --- Three-way rotation: row 136 <- row 137 <- row 138 <- row 136 ---
$saved136 = Get-MatchRowValues $ws 136 $matchCols
$saved137 = Get-MatchRowValues $ws 137 $matchCols
$saved138 = Get-MatchRowValues $ws 138 $matchCols
Set-MatchRowValues $ws 136 $matchCols $saved137
Set-MatchRowValues $ws 137 $matchCols $saved138
Set-MatchRowValues $ws 138 $matchCols $saved136

Write-Host "Row data restored for pairs (7,10) (23,24) (29,30) (72,73) (127,129) and rotation (136,137,138)"
